$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift accelerometer/gyroscope readings: rows 13-22 of the old data become rows 2-11,
# and rows 12-21 receive newly collected "struggle" sensor samples.

$ws.Range("C2").Value = -4.876684755086883
$ws.Range("D2").Value = -2.076402157545074
$ws.Range("E2").Value = -9.344285488128664
$ws.Range("F2").Value = 3.110678434371948
$ws.Range("G2").Value = 1.813658952713013
$ws.Range("H2").Value = -1.52746856212616

$ws.Range("C3").Value = 1.809403181076058
$ws.Range("D3").Value = 13.51332342624669
$ws.Range("E3").Value = -4.600637912750224
$ws.Range("F3").Value = -0.8225300312042236
$ws.Range("G3").Value = 4.032474040985107
$ws.Range("H3").Value = -1.342376351356506

$ws.Range("C4").Value = 1.417896926403043
$ws.Range("D4").Value = 18.05351507663726
$ws.Range("E4").Value = -2.601601481437687
$ws.Range("F4").Value = -0.784503698348999
$ws.Range("G4").Value = -1.800067186355591
$ws.Range("H4").Value = -0.266642689704895

$ws.Range("C5").Value = -2.472576022148132
$ws.Range("D5").Value = 6.250693678855895
$ws.Range("E5").Value = -9.09222507476807
$ws.Range("F5").Value = -0.7644978761672974
$ws.Range("G5").Value = -0.5174028873443604
$ws.Range("H5").Value = 0.9405797719955444

$ws.Range("C6").Value = -2.629937916994093
$ws.Range("D6").Value = 4.76309457421303
$ws.Range("E6").Value = -12.27247226238251
$ws.Range("F6").Value = 0.00137444678694
$ws.Range("G6").Value = -0.2606867551803589
$ws.Range("H6").Value = 3.316082000732422

$ws.Range("C7").Value = -7.906721115112305
$ws.Range("D7").Value = -1.656690835952763
$ws.Range("E7").Value = -8.986274719238274
$ws.Range("F7").Value = 1.413389444351196
$ws.Range("G7").Value = 1.9312504529953
$ws.Range("H7").Value = 2.964529037475586

$ws.Range("C8").Value = -7.878847062587737
$ws.Range("D8").Value = -3.47562119364739
$ws.Range("E8").Value = -5.723703801631919
$ws.Range("F8").Value = 0.1264491081237793
$ws.Range("G8").Value = 3.219565391540528
$ws.Range("H8").Value = -0.1656972020864486

$ws.Range("C9").Value = -3.529585599899285
$ws.Range("D9").Value = -1.653401762247083
$ws.Range("E9").Value = -10.57092201709748
$ws.Range("F9").Value = 1.272585034370422
$ws.Range("G9").Value = 0.1452332139015197
$ws.Range("H9").Value = 1.545947194099426

$ws.Range("C10").Value = -8.618583738803876
$ws.Range("D10").Value = -1.921446576714518
$ws.Range("E10").Value = -7.590942263603199
$ws.Range("F10").Value = 5.68990421295166
$ws.Range("G10").Value = 0.4735732674598694
$ws.Range("H10").Value = 2.645351886749268

$ws.Range("C11").Value = -6.864823818206775
$ws.Range("D11").Value = -3.412744522094729
$ws.Range("E11").Value = -9.644531726837165
$ws.Range("F11").Value = 1.004720568656921
$ws.Range("G11").Value = -1.029460668563843
$ws.Range("H11").Value = 0.8364272117614746

$ws.Range("C12").Value = -0.3932898044586127
$ws.Range("D12").Value = -2.516493201255787
$ws.Range("E12").Value = -11.63133525848387
$ws.Range("F12").Value = -0.5496259927749634
$ws.Range("G12").Value = 1.477530360221863
$ws.Range("H12").Value = -3.358078956604004

$ws.Range("C13").Value = -0.3210607767105248
$ws.Range("D13").Value = 2.201687335968037
$ws.Range("E13").Value = -5.495597958564743
$ws.Range("F13").Value = 1.193783402442932
$ws.Range("G13").Value = 5.348583221435547
$ws.Range("H13").Value = -3.039818286895752

$ws.Range("C14").Value = -3.88015073537826
$ws.Range("D14").Value = 7.319447636604309
$ws.Range("E14").Value = -2.401085853576681
$ws.Range("F14").Value = 0.0236710291355848
$ws.Range("G14").Value = -0.4915938079357147
$ws.Range("H14").Value = -0.5149593949317932

$ws.Range("C15").Value = -2.474413871765147
$ws.Range("D15").Value = 7.422795295715336
$ws.Range("E15").Value = -7.479803562164315
$ws.Range("F15").Value = -0.52702397108078
$ws.Range("G15").Value = 3.396258115768433
$ws.Range("H15").Value = -1.487456917762756

$ws.Range("C16").Value = -4.692895889282228
$ws.Range("D16").Value = 8.487199664115906
$ws.Range("E16").Value = -9.337096989154816
$ws.Range("F16").Value = -0.113315500319004
$ws.Range("G16").Value = 3.309820652008057
$ws.Range("H16").Value = -0.5609270334243774

$ws.Range("C17").Value = -3.639542102813719
$ws.Range("D17").Value = 8.090452075004578
$ws.Range("E17").Value = -3.927529096603386
$ws.Range("F17").Value = 2.999042987823486
$ws.Range("G17").Value = 0.8868235945701599
$ws.Range("H17").Value = -1.359175205230713

$ws.Range("C18").Value = -5.387722790241247
$ws.Range("D18").Value = 3.111244738101949
$ws.Range("E18").Value = 8.269636750221276
$ws.Range("F18").Value = 8.29066276550293
$ws.Range("G18").Value = -1.915215253829956
$ws.Range("H18").Value = 0.5499314665794373

$ws.Range("C19").Value = -1.612288236618022
$ws.Range("D19").Value = 0.4680981636047372
$ws.Range("E19").Value = 12.90354442596435
$ws.Range("F19").Value = 2.119396924972534
$ws.Range("G19").Value = 0.5285511612892151
$ws.Range("H19").Value = 0.0710130855441093

$ws.Range("C20").Value = 0.8442984223365615
$ws.Range("D20").Value = -0.4397069215774612
$ws.Range("E20").Value = 9.103015005588517
$ws.Range("F20").Value = 3.94298243522644
$ws.Range("G20").Value = -0.3740022480487823
$ws.Range("H20").Value = -1.55419385433197

$ws.Range("C21").Value = -3.634706258773804
$ws.Range("D21").Value = -4.129897594451904
$ws.Range("E21").Value = 6.328503251075745
$ws.Range("F21").Value = 3.11648178100586
$ws.Range("G21").Value = -1.741576790809631
$ws.Range("H21").Value = 3.701537847518921

# Remove the now-obsolete last data row (row 22), shrinking the used range to A1:H21
$ws.Rows("22").Delete()
